$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.358.81"
$ws.Range("D3").Value = "2.324.39"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'544.99"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "'130.86"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").Value = "2.321.64"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'23.61"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "2.738.91"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "60.314.40"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "2.320.12"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "'10.58"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "'313.99"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'63.80"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").Value = "'0.994"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "'7.86"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("D29").Value = "'173.29"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("E30").Value = "  +9.32%  "
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").Value = "'5.94"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").Value = "'1.38"
$ws.Range("E34").Value = "  +10.45%  "
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").Value = "'323.00"
$ws.Range("E40").Value = "  +10.76%  "
$ws.Range("D42").Value = "'37.95"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").Value = "'137.93"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").Value = "0.0₆0214"
$ws.Range("E50").Value = "  +18.09%  "
$ws.Range("D51").Value = "'11.03"
$ws.Range("E51").Value = "  +0.68%  "
